$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# The column-type "DATE" was renamed to "TIMESTAMP" for the
# CreateDate (row 21) and LastUpdate (row 23) field definitions.
$ws.Range("D21").Value = "TIMESTAMP"
$ws.Range("D23").Value = "TIMESTAMP"

# Keep the selection where the author left it when the file was saved.
$ws.Range("D23").Select()
